$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: set up the new Rectum cells (AE1:AI1) using the style currently on W1 (Rectum, s=5) ---
$ws.Range("W1").Copy()
$ws.Range("AE1:AI1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AE1").Value = "Rectum24"
$ws.Range("AF1").Value = "Rectum25"
$ws.Range("AG1").Value = "Rectum26"
$ws.Range("AH1").Value = "Rectum27"
$ws.Range("AI1").Value = "Rectum28"

# --- Step 2: turn W1:AD1 into Sigmoid cells (style copied from V1, which is Sigmoid, s=4) ---
$ws.Range("V1").Copy()
$ws.Range("W1:AD1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("W1").Value = "Sigmoid16"
$ws.Range("X1").Value = "Sigmoid17"
$ws.Range("Y1").Value = "Sigmoid18"
$ws.Range("Z1").Value = "Sigmoid19"
$ws.Range("AA1").Value = "Sigmoid20"
$ws.Range("AB1").Value = "Sigmoid21"
$ws.Range("AC1").Value = "Sigmoid22"
$ws.Range("AD1").Value = "Sigmoid23"

# --- Step 3: turn S1:U1 into Descending cells (style copied from R1, which is Descending, s=3) ---
$ws.Range("R1").Copy()
$ws.Range("S1:U1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("S1").Value = "Descending12"
$ws.Range("T1").Value = "Descending13"
$ws.Range("U1").Value = "Descending14"

# --- Step 4: turn P1:Q1 into Transverse cells (style copied from M1, which is Transverse, s=2) ---
$ws.Range("M1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("P1").Value = "Transverse9"
$ws.Range("Q1").Value = "Transverse10"
